$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.180.16'
$ws.Range('E2').Value = '  +2.16%  '
$ws.Range('D3').Value = '2.317.04'
$ws.Range('E3').Value = '  +1.85%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = "'303.24"
$ws.Range('E5').Value = '  +1.66%  '
$ws.Range('D6').Value = "'100.20"
$ws.Range('E6').Value = '  +4.89%  '
$ws.Range('E7').Value = '  +2.80%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  +3.42%  '
$ws.Range('D10').Value = "'34.58"
$ws.Range('E10').Value = '  +3.63%  '
$ws.Range('E11').Value = '  +0.93%  '
$ws.Range('E12').Value = '  +4.20%  '
$ws.Range('D13').Value = "'18.16"
$ws.Range('E13').Value = '  +15.60%  '
$ws.Range('E14').Value = '  +3.39%  '
$ws.Range('D15').Value = '2.689.97'
$ws.Range('E15').Value = '  +2.33%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '2.325.64'
$ws.Range('E16').Value = '  +1.87%  '
$ws.Range('B17').Value = 'Polygon'
$ws.Range('C17').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D17').Value = "'0.822"
$ws.Range('E17').Value = '  +5.24%  '
$ws.Range('D18').Value = '43.132.62'
$ws.Range('E18').Value = '  +2.16%  '
$ws.Range('D19').Value = "'12.59"
$ws.Range('E19').Value = '  +9.23%  '
$ws.Range('D20').Value = '0.0₃0906'
$ws.Range('E20').Value = '  +1.75%  '
$ws.Range('E21').Value = '  +2.21%  '
$ws.Range('D22').Value = "'67.90"
$ws.Range('D23').Value = "'237.66"
$ws.Range('E23').Value = '  +1.77%  '
$ws.Range('E24').Value = '  +12.52%  '
$ws.Range('E25').Value = '  +0.80%  '
$ws.Range('D26').Value = "'1.00"
$ws.Range('E26').Value = '  -0.15%  '
$ws.Range('D27').Value = "'24.89"
$ws.Range('E27').Value = '  +3.99%  '
$ws.Range('D28').Value = "'168.20"
$ws.Range('E28').Value = '  +0.57%  '
$ws.Range('E29').Value = '  -9.18%  '
$ws.Range('D30').Value = "'34.25"
$ws.Range('E30').Value = '  +0.62%  '
$ws.Range('D31').Value = "'9.21"
$ws.Range('E31').Value = '  +1.49%  '
$ws.Range('E32').Value = '  +0.11%  '
$ws.Range('E33').Value = '  +2.73%  '
$ws.Range('D34').Value = "'4.72"
$ws.Range('E34').Value = '  +4.32%  '
$ws.Range('D35').Value = "'2.42"
$ws.Range('E35').Value = '  +4.47%  '
$ws.Range('D36').Value = "'17.18"
$ws.Range('E36').Value = '  +5.27%  '
$ws.Range('D37').Value = "'0.0694"
$ws.Range('E37').Value = '  +0.36%  '
$ws.Range('E38').Value = '  +3.92%  '
$ws.Range('D39').Value = "'1.81"
$ws.Range('E39').Value = '  +4.25%  '
$ws.Range('D40').Value = "'2.83"
$ws.Range('E40').Value = '  +1.22%  '
$ws.Range('E41').Value = '  +0.72%  '
$ws.Range('E42').Value = '  -1.81%  '
$ws.Range('D43').Value = '2.006.96'
$ws.Range('E43').Value = '  +2.25%  '
$ws.Range('E44').Value = '  +3.31%  '
$ws.Range('D45').Value = "'10.16"
$ws.Range('E45').Value = '  +5.62%  '
$ws.Range('D46').Value = "'17.73"
$ws.Range('E46').Value = '  +1.14%  '
$ws.Range('E47').Value = '  +2.44%  '
$ws.Range('D48').Value = "'55.95"
$ws.Range('E48').Value = '  +6.83%  '
$ws.Range('D49').Value = '2.530.41'
$ws.Range('E49').Value = '  +1.21%  '
$ws.Range('E50').Value = '  +4.71%  '
$ws.Range('D51').Value = "'4.59"
$ws.Range('E51').Value = '  +1.34%  '
